# Add three new worksheets to the workbook:
#   - CypherOutput_Message  (copy of the Message sheet content)
#   - StatOutput            (summary counts table)
#   - StatOutput_Message    (Message-sheet content, followed by a second
#                            connection/cypher block for the stats query)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Grab the common "message" values already present on the Message sheet
# so we don't have to retype the connection details.
# ---------------------------------------------------------------------
$msgSheet = $wb.Worksheets.Item("Message")

$neo4jUrlLabel = $msgSheet.Range("A1").Value2
$neo4jUrl      = $msgSheet.Range("A2").Value2
$userNameLabel = $msgSheet.Range("A3").Value2
$userName      = $msgSheet.Range("A4").Value2
$pwdLabel      = $msgSheet.Range("A5").Value2
$pwd           = $msgSheet.Range("A6").Value2
$cypherLabel   = $msgSheet.Range("A7").Value2
$cypherText    = $msgSheet.Range("A8").Value2
$outputLabel   = $msgSheet.Range("A9").Value2
$outputPath    = $msgSheet.Range("A10").Value2

# ---------------------------------------------------------------------
# 1) CypherOutput_Message  -- identical content to the Message sheet
# ---------------------------------------------------------------------
$cypherOutputMessage = $wb.Worksheets.Add($null, $msgSheet)
$cypherOutputMessage.Name = "CypherOutput_Message"

$cypherOutputMessage.Range("A1").Value  = $neo4jUrlLabel
$cypherOutputMessage.Range("A2").Value  = $neo4jUrl
$cypherOutputMessage.Range("A3").Value  = $userNameLabel
$cypherOutputMessage.Range("A4").Value  = $userName
$cypherOutputMessage.Range("A5").Value  = $pwdLabel
$cypherOutputMessage.Range("A6").Value  = $pwd
$cypherOutputMessage.Range("A7").Value  = $cypherLabel
$cypherOutputMessage.Range("A8").Value  = $cypherText
$cypherOutputMessage.Range("A9").Value  = $outputLabel
$cypherOutputMessage.Range("A10").Value = $outputPath

# ---------------------------------------------------------------------
# 2) StatOutput -- summary counts table
# ---------------------------------------------------------------------
$statOutput = $wb.Worksheets.Add($null, $cypherOutputMessage)
$statOutput.Name = "StatOutput"

$statOutput.Range("A1").Value = "number_of_files"
$statOutput.Range("B1").Value = "number_of_sample"
$statOutput.Range("C1").Value = "number_of_cases"
$statOutput.Range("D1").Value = "number_of_study"

$statOutput.Range("A2:D2").NumberFormat = "@"
$statOutput.Range("A2").Value = "1"
$statOutput.Range("B2").Value = "2"
$statOutput.Range("C2").Value = "1"
$statOutput.Range("D2").Value = "1"

# ---------------------------------------------------------------------
# 3) StatOutput_Message -- Message content, then the stats query block
# ---------------------------------------------------------------------
$statCypherText = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Staffordshire Bull Terrier']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

$statOutputMessage = $wb.Worksheets.Add($null, $statOutput)
$statOutputMessage.Name = "StatOutput_Message"

$statOutputMessage.Range("A1").Value  = $neo4jUrlLabel
$statOutputMessage.Range("A2").Value  = $neo4jUrl
$statOutputMessage.Range("A3").Value  = $userNameLabel
$statOutputMessage.Range("A4").Value  = $userName
$statOutputMessage.Range("A5").Value  = $pwdLabel
$statOutputMessage.Range("A6").Value  = $pwd
$statOutputMessage.Range("A7").Value  = $cypherLabel
$statOutputMessage.Range("A8").Value  = $cypherText
$statOutputMessage.Range("A9").Value  = $outputLabel
$statOutputMessage.Range("A10").Value = $outputPath

$statOutputMessage.Range("A11").Value = $neo4jUrlLabel
$statOutputMessage.Range("A12").Value = $neo4jUrl
$statOutputMessage.Range("A13").Value = $userNameLabel
$statOutputMessage.Range("A14").Value = $userName
$statOutputMessage.Range("A15").Value = $pwdLabel
$statOutputMessage.Range("A16").Value = $pwd
$statOutputMessage.Range("A17").Value = $cypherLabel
$statOutputMessage.Range("A18").Value = $statCypherText
$statOutputMessage.Range("A19").Value = $outputLabel
$statOutputMessage.Range("A20").Value = $outputPath

# ---------------------------------------------------------------------
# Sheet ordering ends up as: CypherOutput, Message, CypherOutput_Message,
#                             StatOutput, StatOutput_Message
# (each sheet was inserted immediately After its predecessor above.)
# ---------------------------------------------------------------------
$wb.Worksheets.Item("CypherOutput").Activate()
